# Rename the first worksheet from "Tmed" to "Tmean" and make it the
# active/selected sheet (it was previously the third sheet, "Tmin", that
# was active).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Name = "Tmean"
$ws.Activate()
